$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 8
$ws.Range("C9").Value = 6
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 6
$ws.Range("C12").Value = 4
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 5
$ws.Range("C18").Value = 9
